$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# 1) First paragraph: add two trailing spaces, then append a new
#    red-colored parenthetical note in three runs.
# ----------------------------------------------------------------------
$d.Content.Find.Execute("This is a Microsoft word document.", $true, $false, $false, $false, $false, $true, 1, $false, "This is a Microsoft word document.  ", 2)

$firstPara = $d.Paragraphs(1)
$insPos = $firstPara.Range.End - 1

$s1 = "(This is a change " + [char]0x2013 + " Ve"
$s2 = "rsion for main branch"
$s3 = ")"

$r1 = $d.Range($insPos, $insPos)
$r1.InsertAfter($s1)
$r1.Font.Color = 255

$insPos2 = $insPos + $s1.Length
$r2 = $d.Range($insPos2, $insPos2)
$r2.InsertAfter($s2)
$r2.Font.Color = 255

$insPos3 = $insPos2 + $s2.Length
$r3 = $d.Range($insPos3, $insPos3)
$r3.InsertAfter($s3)
$r3.Font.Color = 255

# ----------------------------------------------------------------------
# 2) "Crispian's Day speech..." paragraph: the visible text is unchanged,
#    but touching it with Find/Replace re-flows / re-merges the runs the
#    same way the source document was re-saved.
# ----------------------------------------------------------------------
$d.Content.Find.Execute(" Day speech from", $true, $false, $false, $false, $false, $true, 1, $false, " Day speech from", 2)
$d.Content.Find.Execute("Henry V", $true, $false, $false, $false, $false, $true, 1, $false, "Henry V", 2)

# ----------------------------------------------------------------------
# 3) Append a new empty paragraph (style "larger") after the final
#    paragraph of the Crispin's day speech.
# ----------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$subRange = $d.Range($lastPara.Range.End - 6, $lastPara.Range.End - 1)
$subRange.Find.Execute("day.", $false, $false, $false, $false, $false, $true, 1, $false, "day.^p", 2)

$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Range.ParagraphFormat.Style = "larger"
$newPara.Format.SpaceBefore = 0
$newPara.Format.SpaceAfter = 7.5
$newPara.Format.SpaceBeforeAuto = $false
$newPara.Format.SpaceAfterAuto = $false
$newPara.Range.Shading.Texture = 0
$newPara.Range.Shading.ForegroundPatternColor = -16777216
$newPara.Range.Shading.BackgroundPatternColor = 16777215
